# Update menu and exports
# Repositions several shapes on slide 1 and fixes a typo in the
# beverages list ("Ice Tea" -> "Iced Tea").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$EMU_PER_POINT = 12700
# Half an EMU (expressed in points) nudges the point value just enough
# to survive the engine's internal float32 rounding and land on the
# exact target EMU instead of one EMU short.
$HALF_EMU = 0.5 / $EMU_PER_POINT

function Set-TopEmu {
    param($Shape, [double]$Emu)
    $Shape.Top = ($Emu / $EMU_PER_POINT) + $HALF_EMU
}

function Set-LeftEmu {
    param($Shape, [double]$Emu)
    $Shape.Left = ($Emu / $EMU_PER_POINT) + $HALF_EMU
}

# Shape 3 (id 91) - round rectangle behind the "DaBeast" logo text
Set-TopEmu $s.Shapes.Item(3) 4401355

# Shape 6 (id 94) - "Breakfast" textbox
Set-LeftEmu $s.Shapes.Item(6) 203126
Set-TopEmu  $s.Shapes.Item(6) 1735751

# Shape 7 (id 95) - "Lunch / Dinner" textbox
Set-LeftEmu $s.Shapes.Item(7) 4035504
Set-TopEmu  $s.Shapes.Item(7) 1773334

# Shape 8 (id 96) - background rectangle (bottom-left column)
Set-TopEmu $s.Shapes.Item(8) 6598074

# Shape 9 (id 97) - straight connector line
Set-TopEmu $s.Shapes.Item(9) 2490439

# Shape 11 (id 99) - "Beverages" textbox; also fix the typo in its text
$sh11 = $s.Shapes.Item(11)
Set-TopEmu $sh11 6598068

$tr11 = $sh11.TextFrame.TextRange
$fullText = $tr11.Text
$oldPhrase = "Sprite, Fanta, Ice Tea, etc."
$newPhrase = "Sprite, Fanta, Iced Tea, etc."
$startIdx = $fullText.IndexOf($oldPhrase)
if ($startIdx -ge 0) {
    $run = $tr11.Characters($startIdx + 1, $oldPhrase.Length)
    $run.Text = $newPhrase
}

# Shape 12 (id 101) - website URL textbox
Set-TopEmu $s.Shapes.Item(12) 9591846

# Shape 14 (group id 4) - "Served All Day" bracket group
Set-TopEmu $s.Shapes.Item(14) 1873759

# Shape 15 (id 19) - QR code picture
Set-TopEmu $s.Shapes.Item(15) 8808892
